# repull data, push all data, mean calculation
# Update the dSF (column F) values for the rows whose data was repulled.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of worksheet row number -> new value for column F ("dSF")
$updates = @{
    2  = -1
    19 = 4
    23 = -2
    27 = -2
    33 = 0
    35 = -2
    36 = 1
    41 = 3
    44 = -3
    50 = -2
    53 = -9
    59 = -2
    62 = 1
    67 = -3
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
